{"js": "// Remove the trailing \"footer\" block of the document: the empty paragraph,\n// the \"Ver no Jupiter...\" paragraph and the \"(c) 2020 ...\" paragraph that\n// immediately follow the \"M\u00c1QUINASEscola PRO-TEC\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"M\u00c1QUINASEscola PRO-TEC\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph 'M\u00c1QUINASEscola PRO-TEC'\");\n}\n\n// The three paragraphs to delete are the ones right after the anchor:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: ... Creative Commons Attribution\"\nconst toDelete = [];\nfor (let offset = 1; offset <= 3 && anchorIndex + offset < items.length; offset++) {\n  toDelete.push(items[anchorIndex + offset]);\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"footer\" block of the document: the empty paragraph,\n# the \"Ver no Jupiter...\" paragraph and the \"(c) 2020 ...\" paragraph that\n# immediately follow the \"M\u00c1QUINASEscola PRO-TEC\" paragraph.\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*M\u00c1QUINASEscola PRO-TEC*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph 'M\u00c1QUINASEscola PRO-TEC'\"\n}\n\n# The three paragraphs to delete are the ones right after the anchor:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) \"\u00a9 2020 . Contact: ... Creative Commons Attribution\"\n# Delete starting from the furthest one so earlier indices stay valid.\nfor ($offset = 3; $offset -ge 1; $offset--) {\n    $idx = $anchorIndex + $offset\n    if ($idx -le $d.Paragraphs.Count) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
